$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) First occurrence of "anatómia" (the underlined term before the
#    "- studuje tvar a stavbu organizmov a ich casti" definition)
#    becomes "morfológia".
# ------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("anatómia", $true, $false, $false, $false, $false, $true, 1, $false, "morfológia", 1) | Out-Null

# ------------------------------------------------------------------
# 2) The "......." run (plus the following "\x96" en-dash run) in the
#    next definition line becomes a single run reading "anatómia - ".
# ------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute(".......–", $true, $false, $false, $false, $false, $true, 1, $false, "anatómia - ", 1) | Out-Null

# ------------------------------------------------------------------
# 3) Move the "_GoBack" bookmark from the end of the document to right
#    after the text we just inserted.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$r3 = $d.Content
$r3.Find.Execute("anatómia - ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bkRange = $d.Range($r3.End, $r3.End)
$d.Bookmarks.Add("_GoBack", $bkRange) | Out-Null
